$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name in the header (merged C3:D3)
$ws.Range("C3").Value = 'KailineLima '

# --- Test case 1 (__init__ / Attribute set to input values.) ---
$ws.Range("E7").Value = 'None '
$ws.Range("F7").Value = 'Rectangle =" Red" 5,  6'
$ws.Range("G7").Value = 'Attributes set '

# --- Test case 2 (__init__ / Exception raised when color is blank) ---
$ws.Range("E8").Value = 'None '
$ws.Range("F8").Value = 'Rectangle ="  " 5,  6'
$ws.Range("G8").Value = 'ValueError'

# --- Test case 3 (__init__ / Exception raised when length is not an integer.) ---
$ws.Range("E9").Value = 'None '
$ws.Range("F9").Value = 'Rectangle = " Red", "five", 6'
$ws.Range("G9").Value = 'ValueError'

# --- Test case 4 (__init__ / Exception raised when width is not an integer.) ---
$ws.Range("E10").Value = 'None '
$ws.Range("F10").Value = 'Rectangle = " Red", 5, " six"'
$ws.Range("G10").Value = 'ValueError'

# --- Test case 5 (__str__ / Returns string formatted appropriately) ---
$ws.Range("E11").Value = 'Rectangle =" Red" 5,  6'
$ws.Range("F11").Value = 'None'
$ws.Range("G11").Value = 'The shape color is red.'

# --- Test case 6 (calculate_area / Returns correct calculated value.) ---
$ws.Range("E12").Value = 'Rectangle =" Red" 5,  6'
$ws.Range("F12").Value = 'None'
$ws.Range("G12").Value = 'This rectangle has four sides with the lengths of 5, 6, 5 and 6 centimeters.'

# --- Test case 7 (calculate_perimeter / Returns correct calculated value.) ---
# Row 13's E/F cells were still using the plain (unwrapped) style, so bring
# them up to the same bold / wrap-text / top-aligned look used elsewhere in
# the table before writing the values.
$ws.Range("E13:F13").Font.Bold = $true
$ws.Range("E13:F13").WrapText = $true
$ws.Range("E13:F13").VerticalAlignment = -4160
$ws.Range("E13").Value = 'Rectangle =" Red" 5,  6'
$ws.Range("F13").Value = 'None'

# Restore the selection to where the editor ended up
$ws.Range("E13").Select() | Out-Null
